# Møller Innkjøp stasjoner - add new NAF station rows + formatting tweaks
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Viking
$ws2 = $wb.Worksheets.Item(2)   # NAF

# ---------------------------------------------------------------------
# Sheet "Viking": postnummer for Vollaveien/Oslo becomes a zero-padded
# text value ("0668") instead of the number 668, and the header row
# becomes bold while the old big-font style on A2 is removed.
# ---------------------------------------------------------------------
$ws1.Range("D2").NumberFormat = "@"
$ws1.Range("D2").Value = "0668"
$ws1.Range("A2").Style = "Normal"
$ws1.Rows.Item(2).RowHeight = 16

$ws1.Range("A1:D1").Font.Bold = $true

$ws1.Columns.Item(1).ColumnWidth = 8.666666666666666

# ---------------------------------------------------------------------
# Sheet "NAF": fill in the previously-empty rows for Namsos, Arendal and
# Steinkjær with address/postcode data, add a "Location" column with a
# map link for Steinkjær, and bold the header row.
# ---------------------------------------------------------------------
$ws2.Range("D2:D4").NumberFormat = "@"
$ws2.Range("D8").NumberFormat = "@"

$ws2.Range("B5").Value = "Pinavegen 2"
$ws2.Range("C5").Value = "Namsos"
$ws2.Range("D5").NumberFormat = "@"
$ws2.Range("D5").Value = "7800"

$ws2.Range("B6").Value = "Åsbieveien 14"
$ws2.Range("C6").Value = "Arendal"
$ws2.Range("D6").NumberFormat = "@"
$ws2.Range("D6").Value = "4848"

$ws2.Range("B7").Value = "Sjøfarstgata 8a"
$ws2.Range("C7").Value = "Steinkjær"
$ws2.Range("D7").NumberFormat = "@"
$ws2.Range("D7").Value = "7714"

$ws2.Range("E1").Value = "Location"
$ws2.Range("E7").Value = "https://maps.app.goo.gl/zKH7TFAFp13YRwvM9"

$ws2.Range("A1:E1").Font.Bold = $true

$ws2.Columns.Item(1).ColumnWidth = 10.330729166666666
$ws2.Columns.Item(3).ColumnWidth = 10.330729166666666

# ---------------------------------------------------------------------
# Selections + active sheet/tab (NAF becomes the active tab, with E7
# selected; Viking keeps A2 selected).
# ---------------------------------------------------------------------
$ws1.Range("A2").Select() | Out-Null
$ws2.Activate()
$ws2.Range("E7").Select() | Out-Null
